$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 19 with the new RSVP'ed event details
$ws.Range("A19").Value = "No"
$ws.Range("C19").Value = "Performance"
$ws.Range("B19").Value = "Shen Yun: 5,000 years of civilisation live on stage"
$ws.Range("D19").Value = "This"

# Move the active selection to A20 (ready for the next entry)
$ws.Range("A20").Select()
